$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the footprint text for the "Switch SP4T" / "PE42440" row (C8)
$ws.Range("C8").Value = "QFN-16-1EP_3x3mm_P0.5mm_EP2.7x2.7mm_ThermalVias"

# Update the active selection to match the saved view state
$ws.Range("J27").Select()
